$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new metric data point as row 66
$ws.Range("A66").Value = "2025-04-29 09:31:43"
$ws.Range("B66").Value = 214
